$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 is the "syntok" dependency row (Name/HomePage/Version/Authors/License/LicenseURL).
# Remove it entirely; Excel shifts rows 35-38 up to become 34-37.
$ws.Rows.Item(34).Delete()
